# Apply updates to the "Historical Figures List" sheet reflecting the
# thumbnail code revamp / file removal commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H3").Value  = "Yes"
$ws.Range("F4").Value  = "No"
$ws.Range("F5").Value  = "No"
$ws.Range("F6").Value  = "No"
$ws.Range("F7").Value  = "No"
$ws.Range("F8").Value  = "No"
$ws.Range("H9").Value  = "Yes"
$ws.Range("F10").Value = "No"
$ws.Range("D13").Value = "Yes"
$ws.Range("E13").Value = "Yes"
$ws.Range("D14").Value = "Yes"
$ws.Range("E14").Value = "Yes"

$wb.Save()
